# Actualiza el Estado de Cuenta: elimina registros anteriores, agrega nuevos
# y refresca los totales / novedades.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Totales de cabecera -----------------------------------------------
$ws.Range("E11").Value = 512460     # VALOR MORA
$ws.Range("F13").Value = 3          # Cant. Periodos

# Las columnas de "Novedad de Ingreso" / "Novedad de Retiro" intercambian
# posicion en el encabezado de la tabla.
$ws.Range("H15").Value = "Novedad de Ingreso"
$ws.Range("I15").Value = "Novedad de Retiro"

# --- Base de datos de trabajadores --------------------------------------
# Se eliminan 4 filas (quedan 9 registros en vez de 13) y se reemplazan los
# datos por el nuevo detalle de mora.
$ws.Range("B17:J20").EntireRow.Delete() | Out-Null

$data = @(
    @("CC", "73187940",   "RICARDO PEREZ PEREZ TOBIAS",        "2505", 56940, 877803),
    @("CC", "73185086",   "RICHARD ACEVEDO CARABALLO",         "2507", 56940, 877803),
    @("CC", "73185086",   "RICHARD ACEVEDO CARABALLO",         "2506", 56940, 877803),
    @("CC", "73185086",   "RICHARD ACEVEDO CARABALLO",         "2505", 56940, 877803),
    @("CC", "73105124",   "JOSUE GUILLERMO VARGAS RODRIGUEZ",  "2505", 56940, 877803),
    @("CC", "1047503132", "ERLIN DAVID ARRIETA NARVAEZ",       "2505", 56940, 1423500),
    @("CC", "1001975088", "CLEIVER BLANCO MADERO",             "2505", 56940, 1423500),
    @("CC", "1007028393", "CARLOS JOSE VALLES PUELLO",         "2507", 56940, 1000000),
    @("CC", "1044934824", "FABER ANDRES VILLAGAS LOPEZ",       "2505", 56940, 828116)
)

$r = 16
foreach ($row in $data) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $r = $r + 1
}
